# fix: memperbaiki struktur tabel invoices dan members
#
# Adds a new "group" column (J) to the members sheet that derives a
# human-readable "golongan"/employment-group label from the existing
# "position" column (I), and corrects a few mis-classified "position"
# values along the way (rows 5, 6 and 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a few mis-entered "position" values (column I) -------------------
$ws.Range("I5").Value  = "p3k"
$ws.Range("I6").Value  = "cpns"
$ws.Range("I21").Value = "p3k"

# --- New "group" column (J) ------------------------------------------------
$ws.Range("J1").Value = "group"

$groups = @{
    2  = "PNS Golongan IV"
    3  = "PNS Golongan IV"
    4  = "PNS Golongan IV"
    5  = "p3k"
    6  = "honor"
    7  = "PNS Golongan III"
    8  = "PNS Golongan III"
    9  = "PNS Golongan IV"
    10 = "PNS Golongan IV"
    11 = "p3k"
    12 = "p3k"
    13 = "p3k"
    14 = "p3k"
    15 = "p3k"
    16 = "p3k"
    17 = "p3k"
    18 = "p3k"
    19 = "p3k"
    20 = "p3k"
    21 = "p3k"
    22 = "honor"
    23 = "honor"
    24 = "honor"
    25 = "honor"
    26 = "honor"
    27 = "honor"
    28 = "honor"
    29 = "honor"
}

foreach ($row in $groups.Keys) {
    $ws.Range("J$row").Value = $groups[$row]
}

# Match the source workbook's formatting: the header and every "honor" cell
# carry the sheet's shared "Text" number format (style index 1), matching
# the format already used throughout column I.
$ws.Range("J1").NumberFormat  = "@"
$ws.Range("J6").NumberFormat  = "@"
$ws.Range("J22:J29").NumberFormat = "@"

# --- Column width / view tidy-up -------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 17

$ws.Range("J26").Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
